$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2, shifting all existing data rows down by one.
$ws.Rows.Item(2).Insert(-4121, 0)

# The inserted row inherits the header's bold formatting; clear it so the
# new row looks like a normal data row (matching the rest of the table).
$ws.Range("A2:R2").ClearFormats()

# Populate the new row 2 with the latest weekly price record.
$ws.Range("A2").Value = 7
$ws.Range("B2").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C2").Value = "Ñuble"
$ws.Range("D2").Value = 44922
$ws.Range("D2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E2").Value = 16
$ws.Range("F2").Value = 100112030
$ws.Range("G2").Value = "Poroto granado"
$ws.Range("H2").Value = "Sin especificar"
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 60
$ws.Range("K2").Value = 32000
$ws.Range("L2").Value = 32000
$ws.Range("M2").Value = 32000
$ws.Range("N2").Value = "$/saco 25 kilos"
$ws.Range("O2").Value = "Región del Maule"
$ws.Range("P2").Value = 1280
$ws.Range("Q2").Value = 25
$ws.Range("R2").Value = "Hortaliza"
